$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 353, pushing existing rows 353:377 down to 354:378
$ws.Range("A353:R353").EntireRow.Insert()

# Fill the new row 353 with the new record's data (same constants as the
# surrounding rows for this market/product, with the new-record-specific values)
$ws.Range("A353").Value = 3
$ws.Range("B353").Value = "Femacal de La Calera"
$ws.Range("C353").Value = "Coquimbo"
$ws.Range("D353").Value = 44714
$ws.Range("D353").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E353").Value = 5
$ws.Range("F353").Value = 100112031
$ws.Range("G353").Value = "Poroto verde"
$ws.Range("H353").Value = "Magnum"
$ws.Range("I353").Value = "Primera"
$ws.Range("J353").Value = 85
$ws.Range("K353").Value = 28000
$ws.Range("L353").Value = 29000
$ws.Range("M353").Value = 28471
$ws.Range("N353").Value = "$/malla 25 kilos"
$ws.Range("O353").Value = "Provincia de Limarí"
$ws.Range("P353").Value = 1139
$ws.Range("Q353").Value = 25
$ws.Range("R353").Value = "Hortaliza"
